$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value2 = 44495
$ws.Range("D4").Value2 = 44420
$ws.Range("J4").Value2 = 120
$ws.Range("K4").Value2 = 13000
$ws.Range("L4").Value2 = 14000
$ws.Range("M4").Value2 = 13500
$ws.Range("P4").Value2 = 338
$ws.Range("D5").Value2 = 44491
$ws.Range("J5").Value2 = 100
$ws.Range("K5").Value2 = 11000
$ws.Range("L5").Value2 = 12000
$ws.Range("M5").Value2 = 11500
$ws.Range("O5").Value2 = 'Provincia del Elquí'
$ws.Range("P5").Value2 = 288
$ws.Range("D6").Value2 = 44490
$ws.Range("J6").Value2 = 100
$ws.Range("K6").Value2 = 11000
$ws.Range("L6").Value2 = 12000
$ws.Range("M6").Value2 = 11500
$ws.Range("P6").Value2 = 288
$ws.Range("D7").Value2 = 44473
$ws.Range("J7").Value2 = 160
$ws.Range("K7").Value2 = 11000
$ws.Range("L7").Value2 = 12000
$ws.Range("M7").Value2 = 11500
$ws.Range("P7").Value2 = 288
$ws.Range("D8").Value2 = 44515
$ws.Range("H8").Value2 = 'Madrigal'
$ws.Range("I8").Value2 = 'Primera'
$ws.Range("K8").Value2 = 11000
$ws.Range("L8").Value2 = 12000
$ws.Range("M8").Value2 = 11500
$ws.Range("P8").Value2 = 288
$ws.Range("D9").Value2 = 44503
$ws.Range("J9").Value2 = 160
$ws.Range("D10").Value2 = 44425
$ws.Range("J10").Value2 = 120
$ws.Range("K10").Value2 = 14000
$ws.Range("L10").Value2 = 15000
$ws.Range("M10").Value2 = 14500
$ws.Range("O10").Value2 = 'Región del Maule'
$ws.Range("P10").Value2 = 362
$ws.Range("D11").Value2 = 44427
$ws.Range("K11").Value2 = 13000
$ws.Range("L11").Value2 = 14000
$ws.Range("M11").Value2 = 13500
$ws.Range("P11").Value2 = 338
$ws.Range("D12").Value2 = 44467
$ws.Range("J12").Value2 = 160
$ws.Range("O12").Value2 = 'Provincia de Limarí'
$ws.Range("D13").Value2 = 44482
$ws.Range("J13").Value2 = 120
$ws.Range("D14").Value2 = 44488
$ws.Range("J14").Value2 = 100
$ws.Range("K14").Value2 = 11000
$ws.Range("L14").Value2 = 12000
$ws.Range("M14").Value2 = 11500
$ws.Range("P14").Value2 = 288
$ws.Range("D15").Value2 = 44446
$ws.Range("D16").Value2 = 44498
$ws.Range("J16").Value2 = 60
$ws.Range("K16").Value2 = 10500
$ws.Range("L16").Value2 = 11000
$ws.Range("M16").Value2 = 10750
$ws.Range("P16").Value2 = 269
$ws.Range("D17").Value2 = 44435
$ws.Range("K17").Value2 = 14000
$ws.Range("L17").Value2 = 15000
$ws.Range("M17").Value2 = 14500
$ws.Range("P17").Value2 = 362
$ws.Range("D18").Value2 = 44508
$ws.Range("J18").Value2 = 160
$ws.Range("D19").Value2 = 44505
$ws.Range("D20").Value2 = 44417
$ws.Range("J20").Value2 = 120
$ws.Range("K20").Value2 = 15000
$ws.Range("L20").Value2 = 16000
$ws.Range("M20").Value2 = 15500
$ws.Range("P20").Value2 = 388
$ws.Range("D21").Value2 = 44487
$ws.Range("J21").Value2 = 100
$ws.Range("O21").Value2 = 'Provincia del Elquí'
$ws.Range("D22").Value2 = 44455
$ws.Range("J22").Value2 = 100
$ws.Range("K22").Value2 = 13000
$ws.Range("L22").Value2 = 14000
$ws.Range("M22").Value2 = 13500
$ws.Range("P22").Value2 = 338
$ws.Range("D23").Value2 = 44510
$ws.Range("D24").Value2 = 44484
$ws.Range("D25").Value2 = 44516
$ws.Range("K25").Value2 = 11000
$ws.Range("L25").Value2 = 12000
$ws.Range("M25").Value2 = 11500
$ws.Range("O25").Value2 = 'Provincia del Elquí'
$ws.Range("P25").Value2 = 288
$ws.Range("D26").Value2 = 44426
$ws.Range("J26").Value2 = 120
$ws.Range("K26").Value2 = 13000
$ws.Range("L26").Value2 = 14000
$ws.Range("M26").Value2 = 13500
$ws.Range("O26").Value2 = 'Región del Maule'
$ws.Range("P26").Value2 = 338
$ws.Range("D27").Value2 = 44512
$ws.Range("J27").Value2 = 120
$ws.Range("D28").Value2 = 44454
$ws.Range("K28").Value2 = 13000
$ws.Range("L28").Value2 = 14000
$ws.Range("M28").Value2 = 13500
$ws.Range("P28").Value2 = 338
$ws.Range("D29").Value2 = 44494
$ws.Range("K29").Value2 = 11000
$ws.Range("L29").Value2 = 12000
$ws.Range("M29").Value2 = 11500
$ws.Range("P29").Value2 = 288
$ws.Range("D30").Value2 = 44489
$ws.Range("J30").Value2 = 120
$ws.Range("D31").Value2 = 44399
$ws.Range("H31").Value2 = 'Española'
$ws.Range("I31").Value2 = 'Segunda'
$ws.Range("K31").Value2 = 15500
$ws.Range("L31").Value2 = 16000
$ws.Range("M31").Value2 = 15750
$ws.Range("P31").Value2 = 394
$ws.Range("D32").Value2 = 44475
$ws.Range("J32").Value2 = 120
$ws.Range("D33").Value2 = 44496
$ws.Range("D34").Value2 = 44453
$ws.Range("K34").Value2 = 12500
$ws.Range("L34").Value2 = 13000
$ws.Range("M34").Value2 = 12750
$ws.Range("P34").Value2 = 319
$ws.Range("D35").Value2 = 44468
$ws.Range("J35").Value2 = 60
$ws.Range("K35").Value2 = 12000
$ws.Range("L35").Value2 = 13000
$ws.Range("M35").Value2 = 12500
$ws.Range("P35").Value2 = 312
